$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.433.36'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '1.636.69'
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.24'
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.533'
$ws.Range("E6").Value = '  +4.46%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -4.60%  '

$ws.Range("E9").Value = '  -2.70%  '

$ws.Range("E10").Value = '  -1.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0886'
$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("D12").Value = '1.870.82'
$ws.Range("E12").Value = '  -1.04%  '

$ws.Range("D13").Value = '1.646.04'
$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.580'
$ws.Range("E14").Value = '  +2.52%  '

$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.13'
$ws.Range("E16").Value = '  -2.31%  '

$ws.Range("D17").Value = '27.423.92'
$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.87'
$ws.Range("E18").Value = '  -2.87%  '

$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("E20").Value = '  -1.31%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").Value = '  -3.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.64'
$ws.Range("E23").Value = '  +3.76%  '

$ws.Range("E24").Value = '  -2.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.05'
$ws.Range("E25").Value = '  +2.14%  '

$ws.Range("E26").Value = '  -2.91%  '

$ws.Range("E27").Value = '  +1.31%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  -3.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").Value = '  -0.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0487'
$ws.Range("E31").Value = '  -2.20%  '

$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("E33").Value = '  +3.55%  '

$ws.Range("D34").Value = '1.412.21'
$ws.Range("E34").Value = '  -2.51%  '

$ws.Range("E35").Value = '  +2.14%  '

$ws.Range("E36").Value = '  -1.61%  '

$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.877'
$ws.Range("E38").Value = '  -4.06%  '

$ws.Range("E39").Value = '  -2.00%  '

$ws.Range("E40").Value = '  -2.36%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.812'
$ws.Range("E42").Value = '  +3.09%  '

$ws.Range("E43").Value = '  +0.31%  '

$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.67'
$ws.Range("E45").Value = '  -2.48%  '

$ws.Range("D46").Value = '1.778.95'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.66'
$ws.Range("E47").Value = '  -3.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.72'
$ws.Range("E48").Value = '  -3.32%  '

$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0991'
$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.71'
$ws.Range("E51").Value = '  -0.91%  '
